# Update daily COVID-19 Valais figures (rows 313-330, sheet "Feuil1").
#
# Only the raw input columns are touched: C (new positive cases), and for a
# few rows the hospital-death / extra-hospital-death split in columns L/M.
# Row 330 was a not-yet-filled-in placeholder day that now gets real data
# (C, E, F, G, L, M). Columns B, H, J and K are formulas ("cumulative
# total = yesterday + today") so they recompute automatically once the
# inputs they depend on change - no need to touch them directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumericValue($cell, $value) {
    # Columns L and M are formatted as Text ("@"). Typing a number straight
    # into such a cell (exactly like doing it by hand in Excel) stores it
    # as text, which is not what the source data has (plain numeric cells).
    # Flip the cell to a numeric format just long enough to commit the
    # value as a real number, then restore the original Text format.
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = $origFormat
}

# Row 313: revised daily case count
$ws.Cells.Item(313, 3).Value = 110

# Row 327: one more hospital death recorded (L) -> K (=L+M) follows
Set-NumericValue $ws.Cells.Item(327, 12) 3

# Row 328: revised daily case count and hospital-death count
$ws.Cells.Item(328, 3).Value = 154
Set-NumericValue $ws.Cells.Item(328, 12) 2

# Row 329: revised daily case count, hospital death and extra-hospital death
$ws.Cells.Item(329, 3).Value = 92
Set-NumericValue $ws.Cells.Item(329, 12) 1
Set-NumericValue $ws.Cells.Item(329, 13) 1

# Row 330: newly filled-in day (previously an empty placeholder row)
$ws.Cells.Item(330, 3).Value = 17
$ws.Cells.Item(330, 5).Value = 11
$ws.Cells.Item(330, 6).Value = 8
$ws.Cells.Item(330, 7).Value = 111
Set-NumericValue $ws.Cells.Item(330, 12) 0
Set-NumericValue $ws.Cells.Item(330, 13) 0
